# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Price cells that look numeric are written with a leading apostrophe
# (Excel's quote-prefix) so they stay plain text, same as the source data
# -- several of these "prices" aren't valid numbers at all (thousands
# separators rendered as extra dots, e.g. 25.741.33) and must round-trip
# byte-for-byte as strings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.741.33'
$ws.Range('E2').Value = '  +5.61%  '
$ws.Range('D3').Value = '1.706.14'
$ws.Range('E3').Value = '  +3.39%  '
$ws.Range('D4').Value = '''1.000'
$ws.Range('E4').Value = '  -0.41%  '
$ws.Range('D5').Value = '''332.09'
$ws.Range('E5').Value = '  +6.76%  '
$ws.Range('D6').Value = '''0.9971'
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('D7').Value = '''0.3684'
$ws.Range('E7').Value = '  +1.40%  '
$ws.Range('D8').Value = '''48.58'
$ws.Range('E8').Value = '  +4.04%  '
$ws.Range('D9').Value = '''0.3303'
$ws.Range('E9').Value = '  +1.90%  '
$ws.Range('D10').Value = '''1.169'
$ws.Range('E10').Value = '  +4.41%  '
$ws.Range('D11').Value = '''0.07344'
$ws.Range('E11').Value = '  +4.82%  '
$ws.Range('D12').Value = '''0.9984'
$ws.Range('E12').Value = '  -0.33%  '
$ws.Range('D13').Value = '''6.202'
$ws.Range('E13').Value = '  +4.63%  '
$ws.Range('D14').Value = '''19.95'
$ws.Range('E14').Value = '  +3.06%  '
$ws.Range('D15').Value = '''6.872'
$ws.Range('E15').Value = '  +4.47%  '
$ws.Range('D16').Value = '1.701.88'
$ws.Range('E16').Value = '  +2.95%  '
$ws.Range('D17').Value = '''0.00001069'
$ws.Range('E17').Value = '  +3.10%  '
$ws.Range('D18').Value = '''0.06616'
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('D19').Value = '''81.03'
$ws.Range('E19').Value = '  +3.87%  '
$ws.Range('D20').Value = '''0.9969'
$ws.Range('E20').Value = '  -0.31%  '
$ws.Range('E21').Value = '  +3.96%  '
$ws.Range('D22').Value = '''6.048'
$ws.Range('E22').Value = '  +2.28%  '
$ws.Range('D23').Value = '''12.99'
$ws.Range('E23').Value = '  +4.43%  '
$ws.Range('D24').Value = '25.714.57'
$ws.Range('E24').Value = '  +5.50%  '
$ws.Range('D25').Value = '''2.454'
$ws.Range('E25').Value = '  -1.08%  '
$ws.Range('D26').Value = '''2.492'
$ws.Range('E26').Value = '  +7.59%  '
$ws.Range('D27').Value = '''149.50'
$ws.Range('E27').Value = '  +1.50%  '
$ws.Range('D28').Value = '''19.20'
$ws.Range('E28').Value = '  +3.84%  '
$ws.Range('D29').Value = '''1.301'
$ws.Range('E29').Value = '  +9.76%  '
$ws.Range('D30').Value = '1.892.76'
$ws.Range('E30').Value = '  +3.12%  '
$ws.Range('D31').Value = '''128.46'
$ws.Range('E31').Value = '  +4.02%  '
$ws.Range('D32').Value = '''4.107'
$ws.Range('E32').Value = '  +0.93%  '
$ws.Range('D33').Value = '''5.958'
$ws.Range('E33').Value = '  +6.11%  '
$ws.Range('D34').Value = '''0.08502'
$ws.Range('E34').Value = '  +1.13%  '
$ws.Range('D35').Value = '''1.676'
$ws.Range('E35').Value = '  +0.91%  '
$ws.Range('D36').Value = '''12.77'
$ws.Range('E36').Value = '  +4.73%  '
$ws.Range('D37').Value = '''5.307'
$ws.Range('E37').Value = '  +3.19%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = '''1.272'
$ws.Range('E38').Value = '  +2.55%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '''0.06230'
$ws.Range('E39').Value = '  +3.89%  '
$ws.Range('D40').Value = '''8.541'
$ws.Range('E40').Value = '  +5.01%  '
$ws.Range('D41').Value = '''0.2121'
$ws.Range('E41').Value = '  +3.37%  '
$ws.Range('D42').Value = '''0.02256'
$ws.Range('E42').Value = '  +2.08%  '
$ws.Range('D43').Value = '''14.56'
$ws.Range('E43').Value = '  +15.82%  '
$ws.Range('D44').Value = '''0.6109'
$ws.Range('E44').Value = '  +4.13%  '
$ws.Range('D45').Value = '''0.9974'
$ws.Range('E45').Value = '  -0.25%  '
$ws.Range('D46').Value = '''3.845'
$ws.Range('E46').Value = '  +2.18%  '
$ws.Range('D47').Value = '''0.5844'
$ws.Range('E47').Value = '  +4.60%  '
$ws.Range('D48').Value = '''126.17'
$ws.Range('E48').Value = '  +3.38%  '
$ws.Range('E49').Value = '  +3.45%  '
$ws.Range('D50').Value = '''0.07217'
$ws.Range('E50').Value = '  +4.86%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '''76.89'
$ws.Range('E51').Value = '  +3.42%  '
